$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 data
$ws.Range("A11").Value = "10"
$ws.Range("B11").Value = 45674
$ws.Range("C11").Value = "Selección ENCIET 202502"
$ws.Range("D11").Value = "Angel Gaibor"

# Match formatting of row 10 (reuse same styles, no new numFmt entries)
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to D12
$ws.Range("D12").Select()
